# Hoan tat Bai tap lon - Lap trinh Web
#
# On the title slide (slide 1), add a "GVHD: ThS. To Oai Hung" line
# above the existing student id / name line that lives in the
# subtitle placeholder shape ("Rectangle 3", shape id 4099).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Rectangle 3") {
        $sh = $cand
        break
    }
}
if ($sh -eq $null) {
    $sh = $s.Shapes.Item(2)
}

$tr = $sh.TextFrame.TextRange
$existing = $tr.Text
$tr.Text = "GVHD: ThS. Tô Oai Hùng" + [char]13 + $existing
